$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16, shifting existing rows 16-40 down to 17-41.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly data point.
$ws.Range("A16").Value = 9
$ws.Range("B16").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44679
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = "Frutos de pepita"
$ws.Range("I16").Value = 100104003
$ws.Range("J16").Value = "Membrillo"
$ws.Range("K16").Value = "Champion"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("Q16").Value = "$/caja 18 kilos granel"
$ws.Range("R16").Value = "Región Metropolitana"
$ws.Range("S16").Value = 556
$ws.Range("T16").Value = 18
